$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(81577, "Pietro Vieira", "Marketing", "Outros", 5, 45085, 5381.77),
    @(17975, "Sr. Ravi Lucca Silveira", "Engenharia", "Outros", 8, 45084, 6766.99),
    @(13015, "Benjamim Rios", "P&D", "Consulta medica", 8, 45089, 7075.69),
    @(86452, "Rebeca Cavalcanti", "TI", "Consulta medica", 8, 45092, 8057.86),
    @(64803, "Bruno Novais", "Operacoes", "Outros", 3, 45102, 4713.33),
    @(69705, "Cecilia Pereira", "TI", "Outros", 1, 45085, 8597.540000000001),
    @(88356, "Maria Julia Sá", "Operacoes", "Consulta medica", 6, 45083, 4358.64),
    @(25357, "Lucas Gabriel Silva", "Atendimento ao Cliente", "Outros", 6, 45091, 7554.76),
    @(25433, "João Guilherme Borges", "Marketing", "Consulta medica", 3, 45092, 5602.82),
    @(80682, "Miguel da Rosa", "Financeiro", "Outros", 3, 45091, 2481.2)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}
